$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country ordering swaps (rows re-sorted by "Casos totales" as data changed) ---

# Butan / Camboya swap (Butan's case count overtook Camboya's)
$ws.Range("A187").Value = "Butan"
$ws.Range("A188").Value = "Camboya"

# Timor Oriental / Santa Lucia swap
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Updated case numbers ---

# India (row 5)
$ws.Range("B5").Value = 6145291
$ws.Range("C5").Value = 2272
$ws.Range("D5").Value = 5101397
$ws.Range("E5").Value = 947543

# Tailandia (row 140)
$ws.Range("B140").Value = 3559
$ws.Range("C140").Value = 14
$ws.Range("D140").Value = 3370
$ws.Range("E140").Value = 130

# Butan (now row 187, after swap above)
$ws.Range("B187").Value = 280
$ws.Range("C187").Value = 7
$ws.Range("D187").Value = 210
$ws.Range("E187").Value = 70

# Camboya (now row 188, after swap above)
$ws.Range("B188").Value = 277
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 275
$ws.Range("E188").Value = 2

# --- Updated "last refreshed" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 06:44"
